# "Update countries & provincias Spain"
#
# The COVID snapshot was refreshed (17:50 -> 18:20). New totals pushed several
# countries' ranking around, so the table's row order changes for a handful of
# neighbouring rows (the country names swap / shift) while almost every data
# row also gets updated case/death counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 3 de Abril de 2020 a las 18:20"

# Estados Unidos (row 4) - numbers refreshed, country stays in place
$ws.Cells.Item(4, 2).Value = 259565
$ws.Cells.Item(4, 3).Value = 14688
$ws.Cells.Item(4, 4).Value = 11960
$ws.Cells.Item(4, 5).Value = 241011
$ws.Cells.Item(4, 7).Value = 524
$ws.Cells.Item(4, 8).Value = 6594

# Italia moves up to row 5 (was España), España drops to row 6 (was Italia)
$ws.Cells.Item(5, 1).Value = "Italia"
$ws.Cells.Item(5, 2).Value = 119827
$ws.Cells.Item(5, 3).Value = 4585
$ws.Cells.Item(5, 4).Value = 19758
$ws.Cells.Item(5, 5).Value = 85388
$ws.Cells.Item(5, 6).Value = 4068
$ws.Cells.Item(5, 7).Value = 766
$ws.Cells.Item(5, 8).Value = 14681

$ws.Cells.Item(6, 1).Value = "España"
$ws.Cells.Item(6, 2).Value = 117710
$ws.Cells.Item(6, 3).Value = 5645
$ws.Cells.Item(6, 4).Value = 30513
$ws.Cells.Item(6, 5).Value = 76262
$ws.Cells.Item(6, 6).Value = 6416
$ws.Cells.Item(6, 7).Value = 587
$ws.Cells.Item(6, 8).Value = 10935

# Minor numeric refreshes, no reordering
$ws.Cells.Item(17, 2).Value = 11464
$ws.Cells.Item(17, 3).Value = 335
$ws.Cells.Item(17, 5).Value = 9274

$ws.Cells.Item(20, 2).Value = 8229
$ws.Cells.Item(20, 3).Value = 185
$ws.Cells.Item(20, 5).Value = 7759
$ws.Cells.Item(20, 7).Value = 19
$ws.Cells.Item(20, 8).Value = 343

$ws.Cells.Item(26, 4).Value = 72
$ws.Cells.Item(26, 5).Value = 3966
$ws.Cells.Item(26, 7).Value = 9
$ws.Cells.Item(26, 8).Value = 53

# Luxemburgo moves up to row 36 (was India), India drops to row 37 (was Luxemburgo)
$ws.Cells.Item(36, 1).Value = "Luxemburgo"
$ws.Cells.Item(36, 2).Value = 2612
$ws.Cells.Item(36, 3).Value = 125
$ws.Cells.Item(36, 4).Value = 174
$ws.Cells.Item(36, 5).Value = 2407
$ws.Cells.Item(36, 6).Value = 33
$ws.Cells.Item(36, 7).Value = 1
$ws.Cells.Item(36, 8).Value = 31

$ws.Cells.Item(37, 1).Value = "India"
$ws.Cells.Item(37, 2).Value = 2567
$ws.Cells.Item(37, 3).Value = 24
$ws.Cells.Item(37, 4).Value = 192
$ws.Cells.Item(37, 5).Value = 2303
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(37, 8).Value = 72

# New entry "Argelia" pushed in at row 52, shifting Colombia..Estonia down one
# row each (Catar lands just above Emiratos Arabes Unidos)
$ws.Cells.Item(52, 1).Value = "Argelia"
$ws.Cells.Item(52, 2).Value = 1171
$ws.Cells.Item(52, 3).Value = 185
$ws.Cells.Item(52, 4).Value = 61
$ws.Cells.Item(52, 5).Value = 1005
$ws.Cells.Item(52, 6).Value = 0
$ws.Cells.Item(52, 7).Value = 19
$ws.Cells.Item(52, 8).Value = 105

$ws.Cells.Item(53, 1).Value = "Colombia"
$ws.Cells.Item(53, 2).Value = 1161
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(53, 4).Value = 55
$ws.Cells.Item(53, 5).Value = 1087
$ws.Cells.Item(53, 6).Value = 50
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 19

$ws.Cells.Item(54, 1).Value = "Singapur"
$ws.Cells.Item(54, 2).Value = 1114
$ws.Cells.Item(54, 3).Value = 65
$ws.Cells.Item(54, 4).Value = 282
$ws.Cells.Item(54, 5).Value = 827
$ws.Cells.Item(54, 6).Value = 24
$ws.Cells.Item(54, 8).Value = 5

$ws.Cells.Item(55, 1).Value = "Croacia"
$ws.Cells.Item(55, 2).Value = 1079
$ws.Cells.Item(55, 3).Value = 68
$ws.Cells.Item(55, 4).Value = 92
$ws.Cells.Item(55, 5).Value = 979
$ws.Cells.Item(55, 6).Value = 39
$ws.Cells.Item(55, 7).Value = 1

$ws.Cells.Item(56, 1).Value = "Catar"
$ws.Cells.Item(56, 2).Value = 1075
$ws.Cells.Item(56, 3).Value = 126
$ws.Cells.Item(56, 4).Value = 93
$ws.Cells.Item(56, 5).Value = 979
$ws.Cells.Item(56, 6).Value = 37
$ws.Cells.Item(56, 8).Value = 3

$ws.Cells.Item(57, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(57, 2).Value = 1024
$ws.Cells.Item(57, 3).Value = 0
$ws.Cells.Item(57, 4).Value = 96
$ws.Cells.Item(57, 5).Value = 920
$ws.Cells.Item(57, 6).Value = 2
$ws.Cells.Item(57, 7).Value = 0
$ws.Cells.Item(57, 8).Value = 8

$ws.Cells.Item(58, 1).Value = "Estonia"
$ws.Cells.Item(58, 2).Value = 961
$ws.Cells.Item(58, 3).Value = 103
$ws.Cells.Item(58, 4).Value = 48
$ws.Cells.Item(58, 5).Value = 901
$ws.Cells.Item(58, 6).Value = 16
$ws.Cells.Item(58, 7).Value = 1
$ws.Cells.Item(58, 8).Value = 12

# New entry "Moldavia" pushed in at row 71, shifting Bosnia y Herzegovina and
# Libano down one row each
$ws.Cells.Item(71, 1).Value = "Moldavia"
$ws.Cells.Item(71, 2).Value = 591
$ws.Cells.Item(71, 3).Value = 86
$ws.Cells.Item(71, 4).Value = 26
$ws.Cells.Item(71, 5).Value = 557
$ws.Cells.Item(71, 6).Value = 65
$ws.Cells.Item(71, 7).Value = 2
$ws.Cells.Item(71, 8).Value = 8

$ws.Cells.Item(72, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(72, 2).Value = 574
$ws.Cells.Item(72, 3).Value = 41
$ws.Cells.Item(72, 4).Value = 27
$ws.Cells.Item(72, 5).Value = 530
$ws.Cells.Item(72, 6).Value = 4

$ws.Cells.Item(73, 1).Value = "Libano"
$ws.Cells.Item(73, 2).Value = 508
$ws.Cells.Item(73, 3).Value = 14
$ws.Cells.Item(73, 4).Value = 50
$ws.Cells.Item(73, 5).Value = 441
$ws.Cells.Item(73, 6).Value = 26
$ws.Cells.Item(73, 7).Value = 1
$ws.Cells.Item(73, 8).Value = 17

# Minor numeric refresh, no reordering
$ws.Cells.Item(109, 2).Value = 156
$ws.Cells.Item(109, 3).Value = 5
$ws.Cells.Item(109, 5).Value = 128

# Nueva Caledonia moves up to row 160 (was Haiti), Haiti drops to row 161
$ws.Cells.Item(160, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(160, 3).Value = 0

$ws.Cells.Item(161, 1).Value = "Haiti"
$ws.Cells.Item(161, 3).Value = 2

# Guinea Ecuatorial moves up to row 163 (was Siria), Siria drops to row 164
$ws.Cells.Item(163, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(163, 3).Value = 1
$ws.Cells.Item(163, 4).Value = 1
$ws.Cells.Item(163, 5).Value = 15
$ws.Cells.Item(163, 8).Value = 0

$ws.Cells.Item(164, 1).Value = "Siria"
$ws.Cells.Item(164, 2).Value = 16
$ws.Cells.Item(164, 4).Value = 0
$ws.Cells.Item(164, 8).Value = 2
